$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134; existing rows 134-193 shift down to 135-194.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new data record.
$ws.Range("A134").Value = 4
$ws.Range("B134").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C134").Value = "Los Lagos"
$ws.Range("D134").Value = 44523
$ws.Range("E134").Value = 10
$ws.Range("F134").Value = 100112040
$ws.Range("G134").Value = "Cilantro"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 150
$ws.Range("K134").Value = 15000
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = 15000
$ws.Range("N134").Value = "`$/caja 36 atados"
$ws.Range("O134").Value = "Región Metropolitana"
$ws.Range("P134").Value = 417
$ws.Range("Q134").Value = 36
$ws.Range("R134").Value = "Hortaliza"

# Keep the date style consistent with the other date cells in column D.
$ws.Range("D134").NumberFormat = $ws.Range("D135").NumberFormat
